$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3732763333333333
$ws.Range("H2").Value = 1.119829
$ws.Range("I2").Value = 0.3554258969843855
$ws.Range("J2").Value = 0.3554258969843855
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.008692112698
$ws.Range("R2").Value = 0.07822901428200001
$ws.Range("S2").Value = 0.003309267622421515
$ws.Range("T2").Value = 0.003309267622421515
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3732763333333333
$ws.Range("H3").Value = 1.119829
$ws.Range("I3").Value = 0.3554258969843855
$ws.Range("J3").Value = 0.3554258969843855
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.04961564137577777
$ws.Range("R3").Value = 0.4465407723819999
$ws.Range("S3").Value = 0.01888970395060778
$ws.Range("T3").Value = 0.01888970395060778
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3732763333333333
$ws.Range("H4").Value = 1.119829
$ws.Range("I4").Value = 0.3554258969843855
$ws.Range("J4").Value = 0.3554258969843855
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 0.8752528716804445
$ws.Range("R4").Value = 7.877275845124
$ws.Range("S4").Value = 0.3332269254113562
$ws.Range("T4").Value = 0.3332269254113562
$ws.Range("G5").Value = 0.668317
$ws.Range("I5").Value = 0.6363574327729865
$ws.Range("J5").Value = 0.6363574327729865
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.015562429662
$ws.Range("R5").Value = 0.140061866958
$ws.Range("S5").Value = 0.005924939815669748
$ws.Range("T5").Value = 0.005924939815669749
$ws.Range("G6").Value = 0.668317
$ws.Range("I6").Value = 0.6363574327729865
$ws.Range("J6").Value = 0.6363574327729865
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("Q6").Value = 0.08883225009533333
$ws.Range("R6").Value = 0.799490250858
$ws.Range("S6").Value = 0.03382028044056282
$ws.Range("T6").Value = 0.03382028044056282
$ws.Range("G7").Value = 0.668317
$ws.Range("I7").Value = 0.6363574327729865
$ws.Range("J7").Value = 0.6363574327729865
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 1.567059899617334
$ws.Range("R7").Value = 14.103539096556
$ws.Range("S7").Value = 0.596612212516754
$ws.Range("T7").Value = 0.596612212516754
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008629333333333334
$ws.Range("H8").Value = 0.025888
$ws.Range("I8").Value = 0.008216670242627913
$ws.Range("J8").Value = 0.008216670242627911
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.000200942656
$ws.Range("R8").Value = 0.001808483904
$ws.Range("S8").Value = 0.00007650303770419251
$ws.Range("T8").Value = 0.00007650303770419251
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008629333333333334
$ws.Range("H9").Value = 0.025888
$ws.Range("I9").Value = 0.008216670242627913
$ws.Range("J9").Value = 0.008216670242627911
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 0.001147005233777778
$ws.Range("R9").Value = 0.010323047104
$ws.Range("S9").Value = 0.0004366886871775371
$ws.Range("T9").Value = 0.0004366886871775371
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008629333333333334
$ws.Range("H10").Value = 0.025888
$ws.Range("I10").Value = 0.008216670242627913
$ws.Range("J10").Value = 0.008216670242627911
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 0.02023393423644445
$ws.Range("R10").Value = 0.182105408128
$ws.Range("S10").Value = 0.007703478517746184
$ws.Range("T10").Value = 0.007703478517746182
